# Apply updated dSF (column F) values for the rows that were repulled/recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    7  = 2
    8  = 0
    9  = -1
    19 = 1
    26 = 2
    29 = 2
    30 = 3
    32 = -2
    39 = 0
    45 = 1
    47 = 3
    59 = 1
    67 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
